{"js": "// Insert a new bold \"No. \" run immediately before the \"${no_surat}\" run,\n// inside the paragraph that currently contains only \"${no_surat}\".\n//\n// The target paragraph ends up reading \"No. ${no_surat}\" as TWO runs:\n//   1) a new run \"No. \" with the same Times New Roman / bold / sz 24 formatting\n//   2) the original, untouched \"${no_surat}\" run\n//\n// We locate the paragraph by its placeholder text and insert a flat-OPC\n// OOXML fragment (a single <w:r>) at a collapsed range right before it \u2014\n// this creates a genuinely separate run (matching the diff) instead of\n// merging the new text into the existing run's <w:t>.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text.indexOf(\"${no_surat}\") !== -1);\nif (!target) {\n  throw new Error(\"Could not find the '${no_surat}' paragraph\");\n}\n\nconst startRange = target.getRange(\"Start\");\n\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r>' +\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:b/>' +\n  '<w:noProof/>' +\n  '<w:sz w:val=\"24\"/>' +\n  '<w:szCs w:val=\"24\"/>' +\n  '<w:lang w:val=\"en-AU\"/>' +\n  '</w:rPr>' +\n  '<w:t xml:space=\"preserve\">No. </w:t>' +\n  '</w:r></w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nstartRange.insertOoxml(flatOpc, \"Start\");\nawait context.sync();\n", "ps1": "# Insert a new bold \"No. \" run immediately before the \"${no_surat}\" run,\n# inside the paragraph that currently contains only \"${no_surat}\".\n#\n# The target paragraph ends up reading \"No. ${no_surat}\" as TWO runs:\n#   1) a new run \"No. \" with the same Times New Roman / bold / sz 24 formatting\n#   2) the original, untouched \"${no_surat}\" run\n#\n# We locate \"${no_surat}\" with Find, then build a brand-new COLLAPSED Range\n# at that match's start position and use Range.InsertXML with a flat-OPC\n# (single <w:r>) fragment. InsertXML replaces the target Range's content \u2014\n# on a genuinely collapsed Range that's a pure insertion \u2014 and, unlike\n# Range.InsertBefore/Range.Text, it creates an actual separate <w:r> element\n# instead of merging the new text into the neighbouring run.\n\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute('${no_surat}')\nif (-not $found) {\n    throw \"Could not find '`${no_surat}' in the document\"\n}\n\n$insertPos = $searchRange.Start\n$target = $d.Range($insertPos, $insertPos)\n\n$flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p><w:r>' +\n    '<w:rPr>' +\n    '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n    '<w:b/>' +\n    '<w:noProof/>' +\n    '<w:sz w:val=\"24\"/>' +\n    '<w:szCs w:val=\"24\"/>' +\n    '<w:lang w:val=\"en-AU\"/>' +\n    '</w:rPr>' +\n    '<w:t xml:space=\"preserve\">No. </w:t>' +\n    '</w:r></w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n$target.InsertXML($flatOpc)\n"}
